# A new weekly price observation was inserted as row 104 (pushing the
# existing rows 104-197 down to 105-198). Replicate that with a native
# row insert followed by populating the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 104, shifting rows 104:197 down to 105:198.
$ws.Rows("104:104").Insert()

# Populate the newly inserted row 104 with the new observation.
$ws.Cells.Item(104, 1).Value = 6
$ws.Cells.Item(104, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(104, 3).Value = "Metropolitana"
$ws.Cells.Item(104, 4).Value = 44658
$ws.Cells.Item(104, 5).Value = 13
$ws.Cells.Item(104, 6).Value = 100112001
$ws.Cells.Item(104, 7).Value = "Berenjena"
$ws.Cells.Item(104, 8).Value = "Sin especificar"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 570
$ws.Cells.Item(104, 11).Value = 8000
$ws.Cells.Item(104, 12).Value = 12000
$ws.Cells.Item(104, 13).Value = 10316
$ws.Cells.Item(104, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(104, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(104, 16).Value = 172
$ws.Cells.Item(104, 17).Value = 60
$ws.Cells.Item(104, 18).Value = "Hortaliza"
